# Updated cryptos list on Wed Feb 21 09:15:22 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.483.58"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "2.921.10"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "362.67"
$ws.Range("E5").Value = "  +1.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.32"
$ws.Range("E6").Value = "  -4.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.546"
$ws.Range("E7").Value = "  -3.25%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.594"
$ws.Range("E9").Value = "  -4.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.21"
$ws.Range("E10").Value = "  -4.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.139"
$ws.Range("E11").Value = "  +1.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0842"
$ws.Range("E12").Value = "  -3.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.66"
$ws.Range("E13").Value = "  -4.16%  "
$ws.Range("D14").Value = "3.390.91"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.40"
$ws.Range("E15").Value = "  -4.76%  "
$ws.Range("D16").Value = "2.942.24"
$ws.Range("E16").Value = "  +0.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.969"
$ws.Range("E17").Value = "  -0.89%  "
$ws.Range("D18").Value = "51.409.42"
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.31"
$ws.Range("E19").Value = "  -1.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.30"
$ws.Range("E20").Value = "  -3.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.15"
$ws.Range("E21").Value = "  -5.52%  "
$ws.Range("D22").Value = "0.0₃0950"
$ws.Range("E22").Value = "  -2.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.82"
$ws.Range("E23").Value = "  -2.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "262.61"
$ws.Range("E24").Value = "  -2.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.71"
$ws.Range("E25").Value = "  -3.77%  "
$ws.Range("E26").Value = "  -5.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.27"
$ws.Range("E27").Value = "  -2.16%  "
$ws.Range("E28").Value = "  -0.24%  "
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.111"
$ws.Range("E29").Value = "  +5.63%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.32"
$ws.Range("E30").Value = "  -5.53%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "10.06"
$ws.Range("E31").Value = "  -4.13%  "
$ws.Range("B32").Value = "RenderToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.21"
$ws.Range("E32").Value = "  +2.52%  "
$ws.Range("E33").Value = "  -2.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "35.10"
$ws.Range("E34").Value = "  -5.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "51.46"
$ws.Range("E35").Value = "  -1.41%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0427"
$ws.Range("E37").Value = "  -3.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.84"
$ws.Range("E38").Value = "  +4.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.18"
$ws.Range("E39").Value = "  -0.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.15"
$ws.Range("E40").Value = "  -5.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.89"
$ws.Range("E41").Value = "  -5.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.114"
$ws.Range("E42").Value = "  -4.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.65"
$ws.Range("E43").Value = "  -1.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "120.08"
$ws.Range("E44").Value = "  +0.95%  "
$ws.Range("E45").Value = "  -1.38%  "
$ws.Range("D46").Value = "2.090.94"
$ws.Range("E46").Value = "  -1.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.23"
$ws.Range("E47").Value = "  -6.20%  "
$ws.Range("E48").Value = "  -7.08%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "3.229.15"
$ws.Range("E49").Value = "  +0.66%  "
$ws.Range("B50").Value = "TheGraph"
$ws.Range("C50").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.238"
$ws.Range("E50").Value = "  -5.02%  "
$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0315"
$ws.Range("E51").Value = "  -7.40%  "
